$wb = $excel.ActiveWorkbook

# Add the two new rows to the "Temp" sheet (columns: Kode table / nama Table)
$ws = $wb.Worksheets.Item("Temp")
$ws.Range("A4").Value = "temp003"
$ws.Range("B4").Value = "temp table search detail room rate"
$ws.Range("A5").Value = "tem004"
$ws.Range("B5").Value = "temp table search detail room rate"

# Update selection to match the new last-used cell
$ws.Range("B5").Select()

# Make "Temp" the active sheet/tab (moves tabSelected + workbook activeTab)
$ws.Activate()
